$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 101
    3  = 297
    4  = 116
    5  = 299
    6  = 272
    7  = 262
    8  = 242
    9  = 263
    10 = 223
    11 = 284
    12 = 252
    13 = 177
    14 = 175
    15 = 263
}

foreach ($row in $values.Keys) {
    $ws.Range("K$row").Value = $values[$row]
}
